$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B (count), C (image), D (word), E (category) for rows 2-33
# Column A (index) and row 1 (headers) are unchanged.

$ws.Cells.Item(2, 2).Value = 102
$ws.Cells.Item(2, 3).Value = "dog/dog012.jpg"
$ws.Cells.Item(2, 4).Value = "rücken"
$ws.Cells.Item(2, 5).Value = "dog"

$ws.Cells.Item(3, 2).Value = 89
$ws.Cells.Item(3, 3).Value = "flower/flower008.jpg"
$ws.Cells.Item(3, 4).Value = "tagen"
$ws.Cells.Item(3, 5).Value = "flower"

$ws.Cells.Item(4, 2).Value = 81
$ws.Cells.Item(4, 3).Value = "dog/dog002.jpg"
$ws.Cells.Item(4, 4).Value = "runden"
$ws.Cells.Item(4, 5).Value = "dog"

$ws.Cells.Item(5, 2).Value = 58
$ws.Cells.Item(5, 3).Value = "flower/flower026.jpg"
$ws.Cells.Item(5, 4).Value = "lehnen"
$ws.Cells.Item(5, 5).Value = "flower"

$ws.Cells.Item(6, 2).Value = 29
$ws.Cells.Item(6, 3).Value = "dog/dog007.jpg"
$ws.Cells.Item(6, 4).Value = "mieten"
$ws.Cells.Item(6, 5).Value = "dog"

$ws.Cells.Item(7, 2).Value = 45
$ws.Cells.Item(7, 3).Value = "flower/flower018.jpg"
$ws.Cells.Item(7, 4).Value = "starten"
$ws.Cells.Item(7, 5).Value = "flower"

$ws.Cells.Item(8, 2).Value = 4
$ws.Cells.Item(8, 3).Value = "flower/flower007.jpg"
$ws.Cells.Item(8, 4).Value = "gelten"
$ws.Cells.Item(8, 5).Value = "flower"

$ws.Cells.Item(9, 2).Value = 82
$ws.Cells.Item(9, 3).Value = "dog/dog019.jpg"
$ws.Cells.Item(9, 4).Value = "fesseln"
$ws.Cells.Item(9, 5).Value = "dog"

$ws.Cells.Item(10, 2).Value = 44
$ws.Cells.Item(10, 3).Value = "dog/dog023.jpg"
$ws.Cells.Item(10, 4).Value = "füttern"
$ws.Cells.Item(10, 5).Value = "dog"

$ws.Cells.Item(11, 2).Value = 13
$ws.Cells.Item(11, 3).Value = "dog/dog031.jpg"
$ws.Cells.Item(11, 4).Value = "stechen"
$ws.Cells.Item(11, 5).Value = "dog"

$ws.Cells.Item(12, 2).Value = 47
$ws.Cells.Item(12, 3).Value = "flower/flower024.jpg"
$ws.Cells.Item(12, 4).Value = "hupen"
$ws.Cells.Item(12, 5).Value = "flower"

$ws.Cells.Item(13, 2).Value = 53
$ws.Cells.Item(13, 3).Value = "flower/flower005.jpg"
$ws.Cells.Item(13, 4).Value = "pflegen"
$ws.Cells.Item(13, 5).Value = "flower"

$ws.Cells.Item(14, 2).Value = 61
$ws.Cells.Item(14, 3).Value = "dog/dog009.jpg"
$ws.Cells.Item(14, 4).Value = "gründen"
$ws.Cells.Item(14, 5).Value = "dog"

$ws.Cells.Item(15, 2).Value = 21
$ws.Cells.Item(15, 3).Value = "dog/dog016.jpg"
$ws.Cells.Item(15, 4).Value = "regnen"
$ws.Cells.Item(15, 5).Value = "dog"

$ws.Cells.Item(16, 2).Value = 98
$ws.Cells.Item(16, 3).Value = "dog/dog018.jpg"
$ws.Cells.Item(16, 4).Value = "saufen"
$ws.Cells.Item(16, 5).Value = "dog"

$ws.Cells.Item(17, 2).Value = 107
$ws.Cells.Item(17, 3).Value = "dog/dog010.jpg"
$ws.Cells.Item(17, 4).Value = "wenden"
$ws.Cells.Item(17, 5).Value = "dog"

$ws.Cells.Item(18, 2).Value = 108
$ws.Cells.Item(18, 3).Value = "dog/dog022.jpg"
$ws.Cells.Item(18, 4).Value = "dauern"
$ws.Cells.Item(18, 5).Value = "dog"

$ws.Cells.Item(19, 2).Value = 92
$ws.Cells.Item(19, 3).Value = "dog/dog015.jpg"
$ws.Cells.Item(19, 4).Value = "füllen"
$ws.Cells.Item(19, 5).Value = "dog"

$ws.Cells.Item(20, 2).Value = 125
$ws.Cells.Item(20, 3).Value = "flower/flower031.jpg"
$ws.Cells.Item(20, 4).Value = "krachen"
$ws.Cells.Item(20, 5).Value = "flower"

$ws.Cells.Item(21, 2).Value = 79
$ws.Cells.Item(21, 3).Value = "flower/flower014.jpg"
$ws.Cells.Item(21, 4).Value = "laufen"
$ws.Cells.Item(21, 5).Value = "flower"

$ws.Cells.Item(22, 2).Value = 28
$ws.Cells.Item(22, 3).Value = "flower/flower015.jpg"
$ws.Cells.Item(22, 4).Value = "wiegen"
$ws.Cells.Item(22, 5).Value = "flower"

$ws.Cells.Item(23, 2).Value = 113
$ws.Cells.Item(23, 3).Value = "flower/flower002.jpg"
$ws.Cells.Item(23, 4).Value = "ehren"
$ws.Cells.Item(23, 5).Value = "flower"

$ws.Cells.Item(24, 2).Value = 122
$ws.Cells.Item(24, 3).Value = "flower/flower020.jpg"
$ws.Cells.Item(24, 4).Value = "bitten"
$ws.Cells.Item(24, 5).Value = "flower"

$ws.Cells.Item(25, 2).Value = 99
$ws.Cells.Item(25, 3).Value = "dog/dog000.jpg"
$ws.Cells.Item(25, 4).Value = "drehen"
$ws.Cells.Item(25, 5).Value = "dog"

$ws.Cells.Item(26, 2).Value = 67
$ws.Cells.Item(26, 3).Value = "flower/flower012.jpg"
$ws.Cells.Item(26, 4).Value = "loben"
$ws.Cells.Item(26, 5).Value = "flower"

$ws.Cells.Item(27, 2).Value = 56
$ws.Cells.Item(27, 3).Value = "dog/dog013.jpg"
$ws.Cells.Item(27, 4).Value = "drohen"
$ws.Cells.Item(27, 5).Value = "dog"

$ws.Cells.Item(28, 2).Value = 85
$ws.Cells.Item(28, 3).Value = "flower/flower011.jpg"
$ws.Cells.Item(28, 4).Value = "backen"
$ws.Cells.Item(28, 5).Value = "flower"

$ws.Cells.Item(29, 2).Value = 39
$ws.Cells.Item(29, 3).Value = "flower/flower022.jpg"
$ws.Cells.Item(29, 4).Value = "jubeln"
$ws.Cells.Item(29, 5).Value = "flower"

$ws.Cells.Item(30, 2).Value = 126
$ws.Cells.Item(30, 3).Value = "dog/dog024.jpg"
$ws.Cells.Item(30, 4).Value = "langen"
$ws.Cells.Item(30, 5).Value = "dog"

$ws.Cells.Item(31, 2).Value = 93
$ws.Cells.Item(31, 3).Value = "flower/flower023.jpg"
$ws.Cells.Item(31, 4).Value = "strahlen"
$ws.Cells.Item(31, 5).Value = "flower"

$ws.Cells.Item(32, 2).Value = 94
$ws.Cells.Item(32, 3).Value = "dog/dog027.jpg"
$ws.Cells.Item(32, 4).Value = "scheitern"
$ws.Cells.Item(32, 5).Value = "dog"

$ws.Cells.Item(33, 2).Value = 3
$ws.Cells.Item(33, 3).Value = "flower/flower017.jpg"
$ws.Cells.Item(33, 4).Value = "biegen"
$ws.Cells.Item(33, 5).Value = "flower"
